$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$lo = $ws1.ListObjects.Item(1)
$ws1.Range("AG1:AK1").EntireColumn.Insert()
try {
    $lo.Resize($ws1.Range("A3:BG86"))
    Write-Output "resize ok"
} catch {
    Write-Output ("err: " + $_.Exception.Message)
}
Write-Output ("Table range: " + $lo.Range.Address())
Write-Output ("ListColumns count: " + $lo.ListColumns.Count)
for ($i=30; $i -le 40; $i++) {
    Write-Output ("Col " + $i + ": " + $lo.ListColumns.Item($i).Name)
}
